# Apply cryptos list refresh (values updated, 3 rows reordered)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be force-written as
# TEXT (matching the source data, which stores these as inline strings) so Excel
# COM does not silently reinterpret them as numbers.
$textForceCells = @("D4", "D5", "D6", "D7", "D9", "D10", "D11", "D12", "D13", "D16", "D17", "D19", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D37", "D41", "D42", "D43", "D44", "D45", "D46", "D47", "D48", "D49", "D50", "D51")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.050.17"
$ws.Range("E2").Value = "  -0.38%  "
$ws.Range("D3").Value = "3.866.95"
$ws.Range("E3").Value = "  -1.97%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "519.65"
$ws.Range("E5").Value = "  +5.32%  "
$ws.Range("D6").Value = "142.01"
$ws.Range("E6").Value = "  -3.56%  "
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  -2.49%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "0.718"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  -4.89%  "
$ws.Range("D11").Value = "0.0000325"
$ws.Range("E11").Value = "  -7.67%  "
$ws.Range("D12").Value = "41.76"
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("D13").Value = "10.35"
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "4.463.82"
$ws.Range("E14").Value = "  -2.27%  "
$ws.Range("D15").Value = "3.879.22"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "21.19"
$ws.Range("E16").Value = "  +6.35%  "
$ws.Range("D17").Value = "14.01"
$ws.Range("E17").Value = "  -2.17%  "
$ws.Range("E18").Value = "  -1.94%  "
$ws.Range("D19").Value = "1.20"
$ws.Range("E19").Value = "  +1.53%  "
$ws.Range("D20").Value = "68.860.39"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").Value = "416.64"
$ws.Range("E21").Value = "  -5.16%  "
$ws.Range("D22").Value = "3.47"
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "13.99"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "87.04"
$ws.Range("E24").Value = "  -2.17%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").Value = "11.97"
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("D26").Value = "4.02"
$ws.Range("E26").Value = "  +5.75%  "
$ws.Range("D27").Value = "10.46"
$ws.Range("E27").Value = "  -5.98%  "
$ws.Range("D28").Value = "35.59"
$ws.Range("E28").Value = "  -4.35%  "
$ws.Range("D29").Value = "13.37"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "676.46"
$ws.Range("E30").Value = "  -3.82%  "
$ws.Range("D31").Value = "6.99"
$ws.Range("E31").Value = "  +15.16%  "
$ws.Range("D32").Value = "2.86"
$ws.Range("E32").Value = "  -1.05%  "
$ws.Range("E33").Value = "  -4.67%  "
$ws.Range("D34").Value = "66.39"
$ws.Range("E34").Value = "  +7.32%  "
$ws.Range("E35").Value = "  -4.59%  "
$ws.Range("D36").Value = "0.0₃0861"
$ws.Range("E36").Value = "  -4.48%  "
$ws.Range("D37").Value = "39.34"
$ws.Range("E37").Value = "  -3.73%  "
$ws.Range("E38").Value = "  +10.18%  "
$ws.Range("E39").Value = "  -3.00%  "
$ws.Range("E40").Value = "  +0.31%  "
$ws.Range("D41").Value = "0.997"
$ws.Range("E41").Value = "  -0.35%  "
$ws.Range("D42").Value = "0.0476"
$ws.Range("E42").Value = "  -2.97%  "
$ws.Range("D43").Value = "3.15"
$ws.Range("E43").Value = "  +4.70%  "
$ws.Range("D44").Value = "2.83"
$ws.Range("E44").Value = "  -2.62%  "
$ws.Range("D45").Value = "3.62"
$ws.Range("E45").Value = "  +6.78%  "
$ws.Range("D46").Value = "0.142"
$ws.Range("E46").Value = "  -1.25%  "
$ws.Range("D47").Value = "0.000284"
$ws.Range("E47").Value = "  +18.28%  "
$ws.Range("D48").Value = "3.05"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("D49").Value = "3.31"
$ws.Range("E49").Value = "  -2.18%  "
$ws.Range("B50").Value = "THORChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D50").Value = "8.78"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "143.31"
$ws.Range("E51").Value = "  -0.89%  "
